$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.964.27"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "1.575.01"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "298.65"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("D7").Value = "0.3734"
$ws.Range("E7").Value = "  -1.14%  "
$ws.Range("D8").Value = "0.3549"
$ws.Range("E8").Value = "  -3.18%  "
$ws.Range("D9").Value = "49.87"
$ws.Range("E9").Value = "  +1.64%  "
$ws.Range("D10").Value = "1.001"
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("D11").Value = "1.207"
$ws.Range("E11").Value = "  -5.26%  "
$ws.Range("D12").Value = "0.07935"
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("D13").Value = "21.69"
$ws.Range("E13").Value = "  -6.23%  "
$ws.Range("D14").Value = "6.413"
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "7.245"
$ws.Range("E15").Value = "  -4.59%  "
$ws.Range("D16").Value = "0.00001212"
$ws.Range("E16").Value = "  -4.57%  "
$ws.Range("D17").Value = "1.582.87"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "91.57"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").Value = "0.06735"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").Value = "17.66"
$ws.Range("E20").Value = "  -4.01%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "6.349"
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("D23").Value = "23.019.48"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").Value = "12.58"
$ws.Range("E24").Value = "  -4.13%  "
$ws.Range("D25").Value = "2.359"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "2.806"
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("D27").Value = "20.51"
$ws.Range("E27").Value = "  -2.95%  "
$ws.Range("D28").Value = "147.18"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").Value = "5.164"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").Value = "131.10"
$ws.Range("E30").Value = "  -1.05%  "
$ws.Range("D31").Value = "2.323"
$ws.Range("E31").Value = "  -3.83%  "
$ws.Range("D32").Value = "6.449"
$ws.Range("E32").Value = "  -7.44%  "
$ws.Range("D33").Value = "1.756.72"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").Value = "0.9253"
$ws.Range("E34").Value = "  -5.58%  "
$ws.Range("D35").Value = "0.07284"
$ws.Range("E35").Value = "  -6.15%  "
$ws.Range("D36").Value = "0.02648"
$ws.Range("E36").Value = "  -5.01%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "0.08717"
$ws.Range("E37").Value = "  -1.77%  "
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2466"
$ws.Range("E38").Value = "  -3.57%  "
$ws.Range("D39").Value = "9.839"
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("D40").Value = "5.929"
$ws.Range("E40").Value = "  -6.11%  "
$ws.Range("D41").Value = "1.336"
$ws.Range("E41").Value = "  -4.70%  "
$ws.Range("D42").Value = "0.6822"
$ws.Range("E42").Value = "  -4.96%  "
$ws.Range("D43").Value = "11.68"
$ws.Range("E43").Value = "  -8.73%  "
$ws.Range("D44").Value = "14.65"
$ws.Range("E44").Value = "  -8.33%  "
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6294"
$ws.Range("E46").Value = "  -5.11%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.959"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "2.224"
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "130.75"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.07835"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").Value = "1.176"
$ws.Range("E51").Value = "  +0.38%  "
